$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-37: price/volume updates only (names and links unchanged)
$ws.Cells.Item(2, 4).Value = '25.790.51'
$ws.Cells.Item(2, 5).Value = '  -3.83%  '
$ws.Cells.Item(3, 4).Value = '1.819.34'
$ws.Cells.Item(3, 5).Value = '  -2.85%  '
$ws.Cells.Item(4, 4).Value = '''1.000'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.25%  '
$ws.Cells.Item(5, 4).Value = '''279.48'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -7.25%  '
$ws.Cells.Item(6, 4).Value = '''1.0000'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.22%  '
$ws.Cells.Item(7, 4).Value = '''0.5102'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -4.12%  '
$ws.Cells.Item(8, 4).Value = '''0.3544'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -5.32%  '
$ws.Cells.Item(9, 4).Value = '''44.36'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -2.65%  '
$ws.Cells.Item(10, 4).Value = '''0.06666'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -7.26%  '
$ws.Cells.Item(11, 4).Value = '''20.08'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -7.15%  '
$ws.Cells.Item(12, 4).Value = '''0.8272'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -6.99%  '
$ws.Cells.Item(13, 4).Value = '''0.07908'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -3.32%  '
$ws.Cells.Item(14, 4).Value = '1.821.71'
$ws.Cells.Item(14, 5).Value = '  -3.03%  '
$ws.Cells.Item(15, 4).Value = '''5.088'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -4.08%  '
$ws.Cells.Item(16, 4).Value = '''88.08'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -5.31%  '
$ws.Cells.Item(17, 4).Value = '''0.9999'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.29%  '
$ws.Cells.Item(18, 4).Value = '''14.13'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -4.85%  '
$ws.Cells.Item(19, 4).Value = '''0.000008030'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -5.82%  '
$ws.Cells.Item(20, 4).Value = '''1.0000'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -0.19%  '
$ws.Cells.Item(21, 4).Value = '25.829.58'
$ws.Cells.Item(21, 5).Value = '  -3.86%  '
$ws.Cells.Item(22, 4).Value = '''4.759'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -4.57%  '
$ws.Cells.Item(23, 4).Value = '''9.981'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -6.15%  '
$ws.Cells.Item(24, 4).Value = '''6.126'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -4.16%  '
$ws.Cells.Item(25, 4).Value = '''2.232'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -2.60%  '
$ws.Cells.Item(26, 4).Value = '''142.26'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -3.01%  '
$ws.Cells.Item(27, 4).Value = '''1.669'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -4.17%  '
$ws.Cells.Item(28, 4).Value = '''17.16'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -4.93%  '
$ws.Cells.Item(29, 4).Value = '''109.61'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -3.87%  '
$ws.Cells.Item(30, 4).Value = '''4.323'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -8.18%  '
$ws.Cells.Item(31, 4).Value = '''4.245'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -7.99%  '
$ws.Cells.Item(32, 4).Value = '''0.08757'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -3.86%  '
$ws.Cells.Item(33, 4).Value = '''0.04913'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -2.09%  '
$ws.Cells.Item(34, 4).Value = '''0.7313'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -9.78%  '
$ws.Cells.Item(35, 4).Value = '''1.141'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -2.60%  '
$ws.Cells.Item(36, 4).Value = '''2.873'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -3.13%  '
$ws.Cells.Item(37, 4).Value = '''3.162'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -1.25%  '

# Rows 38-51: Frax inserted at row 38, shifting subsequent coins down by one; Cronos drops off the bottom
$ws.Cells.Item(38, 2).Value = 'Frax'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(38, 4).Value = '''0.9990'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -0.37%  '
$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(39, 4).Value = '''2.367'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -11.08%  '
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(40, 4).Value = '''0.01861'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -4.95%  '
$ws.Cells.Item(41, 2).Value = 'TheSandbox'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(41, 4).Value = '''0.5157'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -15.29%  '
$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(42, 4).Value = '''0.9693'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -9.24%  '
$ws.Cells.Item(43, 2).Value = 'Quant'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(43, 4).Value = '''113.20'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -1.65%  '
$ws.Cells.Item(44, 2).Value = 'FraxShare'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(44, 4).Value = '''6.256'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -4.46%  '
$ws.Cells.Item(45, 2).Value = 'Aptos'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(45, 4).Value = '''8.066'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -8.94%  '
$ws.Cells.Item(46, 2).Value = 'PaxDollar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(46, 4).Value = '''0.9998'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -0.24%  '
$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).Value = '''0.4564'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -11.71%  '
$ws.Cells.Item(48, 2).Value = 'Algorand'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(48, 4).Value = '''0.1370'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -8.42%  '
$ws.Cells.Item(49, 2).Value = 'Elrond'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Cells.Item(49, 4).Value = '''36.63'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -2.36%  '
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value = '''9.234'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -6.87%  '
$ws.Cells.Item(51, 2).Value = 'NEARProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(51, 4).Value = '''1.500'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -8.63%  '
